$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "20÷6="  # was 47÷6=
$t.Cell(1, 2).Range.Text = "81÷2="  # was 37÷8=
$t.Cell(1, 3).Range.Text = "23÷6="  # was 98÷2=
$t.Cell(1, 4).Range.Text = "68÷3="  # was 19÷6=
$t.Cell(1, 5).Range.Text = "69÷5="  # was 12÷8=
$t.Cell(5, 1).Range.Text = "60÷8="  # was 81÷9=
$t.Cell(5, 2).Range.Text = "16÷7="  # was 33÷5=
$t.Cell(5, 3).Range.Text = "85÷8="  # was 94÷4=
$t.Cell(5, 4).Range.Text = "15÷6="  # was 98÷7=
$t.Cell(5, 5).Range.Text = "60÷8="  # was 33÷3=
$t.Cell(9, 1).Range.Text = "74÷4="  # was 85÷8=
$t.Cell(9, 2).Range.Text = "84÷5="  # was 31÷4=
$t.Cell(9, 3).Range.Text = "12÷2="  # was 49÷8=
$t.Cell(9, 4).Range.Text = "78÷3="  # was 97÷5=
$t.Cell(9, 5).Range.Text = "77÷9="  # was 40÷6=
$t.Cell(13, 1).Range.Text = "13÷2="  # was 50÷4=
$t.Cell(13, 2).Range.Text = "32÷4="  # was 22÷9=
$t.Cell(13, 3).Range.Text = "14÷4="  # was 83÷7=
$t.Cell(13, 4).Range.Text = "18÷7="  # was 27÷6=
$t.Cell(13, 5).Range.Text = "78÷7="  # was 54÷4=
$t.Cell(17, 1).Range.Text = "71÷2="  # was 76÷8=
$t.Cell(17, 2).Range.Text = "97÷8="  # was 83÷7=
$t.Cell(17, 3).Range.Text = "71÷8="  # was 27÷3=
$t.Cell(17, 4).Range.Text = "50÷2="  # was 78÷4=
$t.Cell(17, 5).Range.Text = "23÷4="  # was 95÷5=
